$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New combined card text for rows 2-8 (one row per card, combining the
# card name and its attribute list into a single Python-tuple-like string).
$ws.Range("A2").Value = "('Carnophage', ['{B}', 'Creature " + [char]0x2014 + " Zombie', 'At the beginning of your upkeep, tap Carnophage unless you pay 1 life.', '2/2'])"
$ws.Range("A3").Value = "('Fireblast', ['{4}{R}{R}', 'Instant', 'You may sacrifice two Mountains rather than pay this spell" + [char]0x2019 + "s mana cost.', 'Fireblast deals 4 damage to any target.'])"
$ws.Range("A4").Value = "('Impulse', ['{1}{U}', 'Instant', 'Look at the top four cards of your library. Put one of them into your hand and the rest on the bottom of your library in any order.'])"
$ws.Range("A5").Value = "('Jackal Pup', ['{R}', 'Creature " + [char]0x2014 + " Jackal', 'Whenever Jackal Pup is dealt damage, it deals that much damage to you.', '2/1'])"
$ws.Range("A6").Value = "('Ophidian', ['{2}{U}', 'Creature " + [char]0x2014 + " Snake', 'Whenever Ophidian attacks and isn" + [char]0x2019 + "t blocked, you may draw a card. If you do, Ophidian assigns no combat damage this turn.', '1/3'])"
$ws.Range("A7").Value = "('Quirion Ranger', ['{G}', 'Creature " + [char]0x2014 + " Elf', 'Return a Forest you control to its owner" + [char]0x2019 + "s hand: Untap target creature. Activate this ability only once each turn.', '1/1'])"
$ws.Range("A8").Value = "('Swords to Plowshares', ['{W}', 'Instant', 'Exile target creature. Its controller gains life equal to its power.'])"

# Remove the now-unused rows 9-34 so the sheet's used range shrinks to A1:A8.
$ws.Range("A9:A34").ClearContents()
